$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Observações" note in D3: disciplina preferences override curso preferences
$ws.Range("D3").Value = "Preferências de disciplina sobrescrevem preferências de curso"

# New rows of disciplina -> sala preferences (rows 15-21)
$ws.Range("A15").Value = "GCB295"
$ws.Range("B15").Value = "109-DE"

$ws.Range("A16").Value = "GCH632"
$ws.Range("B16").Value = "105-B"

$ws.Range("A17").Value = "GCH626"
$ws.Range("B17").Value = "105-B"

$ws.Range("A18").Value = "GEX557"
$ws.Range("B18").Value = "109-DE"

$ws.Range("A19").Value = "GCH633"
$ws.Range("B19").Value = "109-DE"

$ws.Range("A20").Value = "GCH627"
$ws.Range("B20").Value = "105-B"

$ws.Range("A21").Value = "GEX556"
$ws.Range("B21").Value = "110-DE"

# Highlighted rows (15, 16, 20, 21) get a slightly larger black Arial font
$ws.Range("A15:A16").Font.Color = 0
$ws.Range("A15:A16").Font.Size = 11

$ws.Range("A20:A21").Font.Color = 0
$ws.Range("A20:A21").Font.Size = 11

# Sheet view: scroll back to A1, select A20, zoom to 75%
$ws.Range("A1").Select()
$excel.ActiveWindow.Zoom = 75
$ws.Range("A20").Select()
